# Add new gear code 838 ("Sill-/strömmingryssja/push up, vittjanpåse")
# This inserts a new row at position 49 in Sheet1, pushing the existing
# rows 49-61 down to 50-62, and fills in the new row's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 49 (shifts rows 49:61 down to 50:62)
$ws.Rows.Item(49).EntireRow.Insert() | Out-Null

# Populate the new row with the new gear code data
$ws.Cells.Item(49, 1).Value = "FYK_C"
$ws.Cells.Item(49, 2).Value = "Sill-/strömmingryssja/push up, vittjanpåse"
$ws.Cells.Item(49, 3).Value = 838
$ws.Cells.Item(49, 4).Value = "gd"

# Match formatting used for similarly long wrapped descriptions in this column
$ws.Cells.Item(49, 2).WrapText = $true
$ws.Rows.Item(49).RowHeight = 35.05

# Refresh the autofilter range to include the new last row (A1:D62)
$ws.AutoFilterMode = $false
$ws.Range("A1:D62").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$D`$62"
    }
}

# Restore the active cell selection as recorded after the edit
$ws.Range("B54").Select() | Out-Null
